$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6325.7144
$ws.Range("I51").Value = 4472.25
$ws.Range("K51").Value = 4472.25
$ws.Range("M51").Value = -3988.25
$ws.Range("H106").Value = 1730.8182
$ws.Range("I106").Value = 1730.8182
$ws.Range("K106").Value = 1730.8182
$ws.Range("M106").Value = -1099.8182
$ws.Range("H112").Value = 2446.9
$ws.Range("I112").Value = 679.5
$ws.Range("J112").Value = 2643.2778
$ws.Range("K112").Value = 2038.5
$ws.Range("L112").Value = 7929.8334
$ws.Range("M112").Value = -930.5
$ws.Range("N112").Value = -10145.8334
$ws.Range("H138").Value = 2415.7322
$ws.Range("I138").Value = 1735.5714
$ws.Range("J138").Value = 3095.8928
$ws.Range("K138").Value = 5206.7142
$ws.Range("L138").Value = 9287.678400000001
$ws.Range("M138").Value = -66.71420000000035
$ws.Range("N138").Value = -19567.6784
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1586.75
$ws.Range("I2").Value = 1174
$ws.Range("K2").Value = 1174
$ws.Range("M2").Value = -1061
$ws.Range("H45").Value = 2850.5789
$ws.Range("I45").Value = 1636.4615
$ws.Range("J45").Value = 5481.1665
$ws.Range("K45").Value = 1636.4615
$ws.Range("L45").Value = 5481.1665
$ws.Range("M45").Value = -1259.4615
$ws.Range("N45").Value = -6235.1665
$ws.Range("H61").Value = 7891.3335
$ws.Range("I61").Value = 746.3333
$ws.Range("K61").Value = 746.3333
$ws.Range("M61").Value = -534.3333
$ws.Range("H110").Value = 565.6667
$ws.Range("I110").Value = 576.7778
$ws.Range("K110").Value = 576.7778
$ws.Range("M110").Value = 1468.2222
$ws.Range("H116").Value = 1586.75
$ws.Range("I116").Value = 1174
$ws.Range("K116").Value = 1174
$ws.Range("M116").Value = 1120
$ws.Range("H136").Value = 7891.3335
$ws.Range("I136").Value = 746.3333
$ws.Range("K136").Value = 2238.9999
$ws.Range("M136").Value = 311.0001000000002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1586.75
$ws.Range("I3").Value = 1174
$ws.Range("K3").Value = 1174
$ws.Range("M3").Value = -1060
$ws.Range("H82").Value = 23891.889
$ws.Range("I82").Value = 15005.4
$ws.Range("K82").Value = 15005.4
$ws.Range("M82").Value = -14622.4
$ws.Range("H85").Value = 23891.889
$ws.Range("I85").Value = 15005.4
$ws.Range("K85").Value = 15005.4
$ws.Range("M85").Value = -13679.4
$ws.Range("H107").Value = 5407.0454
$ws.Range("I107").Value = 5401.5264
$ws.Range("K107").Value = 5401.5264
$ws.Range("M107").Value = -3481.5264
$ws.Range("H134").Value = 7729.2856
$ws.Range("I134").Value = 7990.162
$ws.Range("K134").Value = 23970.486
$ws.Range("M134").Value = -21435.486
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1963.3077
$ws.Range("I58").Value = 1557.2354
$ws.Range("K58").Value = 1557.2354
$ws.Range("M58").Value = -1354.2354
$ws.Range("H122").Value = 3098.6
$ws.Range("I122").Value = 3098.6
$ws.Range("K122").Value = 9295.799999999999
$ws.Range("M122").Value = -6845.799999999999
$ws.Range("H125").Value = 70599.39999999999
$ws.Range("J125").Value = 70599.39999999999
$ws.Range("L125").Value = 70599.39999999999
$ws.Range("N125").Value = -75519.39999999999
$ws.Range("H136").Value = 1963.3077
$ws.Range("I136").Value = 1557.2354
$ws.Range("K136").Value = 4671.706200000001
$ws.Range("M136").Value = -2121.706200000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 161.19048
$ws.Range("I2").Value = 164.47058
$ws.Range("J2").Value = 147.25
$ws.Range("K2").Value = 986.82348
$ws.Range("L2").Value = 883.5
$ws.Range("M2").Value = -873.82348
$ws.Range("N2").Value = -1109.5
$ws.Range("H37").Value = 41972.555
$ws.Range("J37").Value = 41972.555
$ws.Range("L37").Value = 125917.665
$ws.Range("N37").Value = -126141.665
$ws.Range("H75").Value = 4930.391
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4930.391
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 14791.173
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -16787.173
$ws.Range("H78").Value = 4930.391
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4930.391
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 44373.519
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -54357.519
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11995.143
$ws.Range("I70").Value = 11853.4
$ws.Range("J70").Value = 12349.5
$ws.Range("K70").Value = 11853.4
$ws.Range("L70").Value = 12349.5
$ws.Range("M70").Value = -11583.4
$ws.Range("N70").Value = -12889.5
$ws.Range("H73").Value = 11995.143
$ws.Range("I73").Value = 11853.4
$ws.Range("J73").Value = 12349.5
$ws.Range("K73").Value = 11853.4
$ws.Range("L73").Value = 12349.5
$ws.Range("M73").Value = -10917.4
$ws.Range("N73").Value = -14221.5
$ws.Range("H102").Value = 2295.2083
$ws.Range("I102").Value = 1909.5714
$ws.Range("K102").Value = 1909.5714
$ws.Range("M102").Value = -287.5714
$ws.Range("H122").Value = 3819.4614
$ws.Range("I122").Value = 3993.1
$ws.Range("K122").Value = 11979.3
$ws.Range("M122").Value = -9529.299999999999
$ws.Range("H132").Value = 1995.3793
$ws.Range("I132").Value = 1610.2693
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 4830.8079
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -2300.8079
$ws.Range("N132").Value = -21059
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5333.3335
$ws.Range("I7").Value = 5333.3335
$ws.Range("K7").Value = 5333.3335
$ws.Range("M7").Value = -5221.3335
$ws.Range("H34").Value = 116674.664
$ws.Range("J34").Value = 100012
$ws.Range("L34").Value = 100012
$ws.Range("N34").Value = -100356
$ws.Range("H126").Value = 5333.3335
$ws.Range("I126").Value = 5333.3335
$ws.Range("K126").Value = 16000.0005
$ws.Range("M126").Value = -13530.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 73534.55
$ws.Range("I122").Value = 111412.305
$ws.Range("K122").Value = 334236.915
$ws.Range("M122").Value = -331786.915
$ws.Range("H136").Value = 27950.818
$ws.Range("I136").Value = 34495.883
$ws.Range("K136").Value = 103487.649
$ws.Range("M136").Value = -100937.649
